$d = $word.ActiveDocument

# --- Edit 1: Remove the duplicate "Reveal.js and Spectacle.js" paragraph
# that sits right before the "Mathematics / Operations Research / Network
# Science" heading (the one filed under Engineering / General Management -
# it duplicates the legitimate one under Training / Mentorship).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "Reveal\.js" -and $t -match "Spectacle\.js") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -match "Mathematics / Operations Research / Network Science") {
            [void]$p.Range.Delete()
            break
        }
    }
}

# --- Edit 2: Append the rank/paygrade to the Navy Selected Ready Reserve bullet.
[void]$d.Content.Find.Execute(
    "2014 - Present Navy Selected Ready Reserve member",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2014 - Present Navy Selected Ready Reserve member (LCDR, O-4)", 2
)
